$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.879.43"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "3.388.30"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").Value = "3.380.66"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "3.928.56"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "3.372.25"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "65.609.59"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "463.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "89.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "580.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "62.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.143"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.378"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("D41").Value = "0.0₃0750"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").Value = "3.104.65"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0418"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.38%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.16%  "
